$wb = $excel.ActiveWorkbook

$advisor = "ALMEIDA CUATIN JHONATHANN CARLOS"
$newName = "LLERENA CONDO SANDRA MARISOL"

# --- Sheet "VENTAS POR GRUPO": columns A:R, insert new row 19 ---
$ws1 = $wb.Worksheets.Item("VENTAS POR GRUPO")
$ws1.Rows.Item(19).Insert()

$ws1.Cells.Item(19, 1).Value = $advisor
$ws1.Cells.Item(19, 2).Value = $newName
for ($col = 3; $col -le 18; $col++) {
    $ws1.Cells.Item(19, $col).Value = 0
    $ws1.Cells.Item(19, $col).NumberFormat = '"$"#,##0.00'
}

# Update the "x de y" summary row (now row 35) counts from 32 to 33
for ($col = 3; $col -le 18; $col++) {
    $cell = $ws1.Cells.Item(35, $col)
    $text = [string]$cell.Value2
    $cell.Value = $text.Replace("de 32", "de 33")
}

# --- Sheet "VENTA MENSUAL": columns A:G, insert new row 19 ---
$ws2 = $wb.Worksheets.Item("VENTA MENSUAL")
$ws2.Rows.Item(19).Insert()

$ws2.Cells.Item(19, 1).Value = $advisor
$ws2.Cells.Item(19, 2).Value = $newName
for ($col = 3; $col -le 7; $col++) {
    $ws2.Cells.Item(19, $col).Value = 0
    $ws2.Cells.Item(19, $col).NumberFormat = '"$"#,##0.00'
}
